# Update attendance / price figures for three worksheets:
#   展览 (Exhibitions), 演出 (Performances), 全部类型 (All types)
# as recorded by the new site generation run.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 457
$wsExpo.Range("F3").Value = 5515
$wsExpo.Range("G3").Value = 62
$wsExpo.Range("F9").Value = 523

# --- 演出 sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 8

# --- 全部类型 sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 457
$wsAll.Range("F3").Value = 5515
$wsAll.Range("G3").Value = 62
$wsAll.Range("F9").Value = 8
$wsAll.Range("F11").Value = 523
